# EditCommandUML.pptx - "resolved discrepancies in UML diagram for EditStudentCommand"
#
# 1) Refresh the "last modified" date field baked into the slide master and
#    every slide layout (24-10-2018 -> 10-11-2018).
# 2) On slide 1, introduce a new abstract-superclass box ("{abstract}" /
#    "Command") sized to sit above "EditStudentCommand", styled/positioned
#    like the old small "Command" label box, and remove that now-redundant
#    "Command" label box.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Step 1: update the datetimeFigureOut placeholder text wherever it is
# still showing the old date, across the slide master and all layouts.
# ---------------------------------------------------------------------
$oldDate = "24-10-2018"
$newDate = "10-11-2018"

function Update-DateShapes($shapes) {
    for ($j = 1; $j -le $shapes.Count; $j++) {
        $sh = $shapes.Item($j)
        if ($sh.Name -like "Date Placeholder*") {
            if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
                if ($sh.TextFrame.TextRange.Text -eq $oldDate) {
                    $sh.TextFrame.TextRange.Text = $newDate
                }
            }
        }
    }
}

$master = $p.SlideMaster
Update-DateShapes $master.Shapes
for ($i = 1; $i -le $master.CustomLayouts.Count; $i++) {
    $layout = $master.CustomLayouts.Item($i)
    Update-DateShapes $layout.Shapes
}

# ---------------------------------------------------------------------
# Step 2: rework the UML diagram on slide 1.
# ---------------------------------------------------------------------
$slide = $p.Slides.Item(1)

# Locate the two shapes we care about by their stable shape Id (not
# positional index, since index shifts as shapes are added/removed):
#   id 7  -> "EditStudentCommand" box (stays put, becomes the anchor we
#             insert the new box in front of)
#   id 10 -> the small "Command" label box being replaced
$editStudentCommandShape = $null
$oldCommandShape = $null
for ($k = 1; $k -le $slide.Shapes.Count; $k++) {
    $candidate = $slide.Shapes.Item($k)
    if ($candidate.Id -eq 7) { $editStudentCommandShape = $candidate }
    if ($candidate.Id -eq 10) { $oldCommandShape = $candidate }
}

# Target geometry (EMU) for the new box.
$newLeftEmu   = 5361112
$newTopEmu    = 2063931
$newWidthEmu  = 1589103
$newHeightEmu = 590268

# PowerPoint's Shape.Left/Top/Width/Height are expressed in points
# (1 pt = 12700 EMU) and stored as single-precision floats, which loses
# a touch of precision versus the raw EMU value. A tiny nudge keeps the
# round-tripped value snapped to the exact target EMU instead of one
# unit short.
function ConvertTo-Points($emu) {
    return ($emu / 12700) + 0.00001
}

# Duplicate the old "Command" box so the new shape inherits its exact
# style (accent4 line/fill/effect refs) and text formatting, then
# retarget its geometry/content.
$dup = $oldCommandShape.Duplicate()
$newShape = $dup.Item(1)

# Turn the single "Command" paragraph into two paragraphs: "{abstract}"
# followed by "Command", both sharing the same run formatting.
$null = $newShape.TextFrame.TextRange.InsertBefore("{abstract}" + [char]13)

$newShape.Left   = ConvertTo-Points $newLeftEmu
$newShape.Top    = ConvertTo-Points $newTopEmu
$newShape.Width  = ConvertTo-Points $newWidthEmu
$newShape.Height = ConvertTo-Points $newHeightEmu

# Move the new shape to the back of the z-order / front of the XML shape
# list, so it lands immediately before the "EditStudentCommand" shape,
# matching the diff's shape ordering.
$newShape.ZOrder(1)  # msoSendToBack

# Remove the old "Command" label box now that its content lives in the
# new shape.
$oldCommandShape.Delete()
